# Scheduled market-price refresh: update cached average-price / leve-profit
# columns (H, I, J, K, L, M, N) for the affected leve rows on each job sheet.
# Values come from the latest Universalis price pull; formulas are not used
# on these sheets, so every touched cell is a plain literal.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4716

$ws.Range("H43").Value = 1321.6666
$ws.Range("I43").Value = 1321.6666
$ws.Range("K43").Value = 1321.6666
$ws.Range("M43").Value = -1252.6666

$ws.Range("H62").Value = 7953.625
$ws.Range("I62").Value = 7736.8
$ws.Range("K62").Value = 7736.8
$ws.Range("M62").Value = -7112.8

$ws.Range("H65").Value = 7953.625
$ws.Range("I65").Value = 7736.8
$ws.Range("K65").Value = 38684
$ws.Range("M65").Value = -35564

$ws.Range("H86").Value = 20045.238
$ws.Range("I86").Value = 4710.6665
$ws.Range("J86").Value = 31546.166
$ws.Range("K86").Value = 4710.6665
$ws.Range("L86").Value = 31546.166
$ws.Range("M86").Value = -3587.6665
$ws.Range("N86").Value = -33792.166

$ws.Range("H89").Value = 20045.238
$ws.Range("I89").Value = 4710.6665
$ws.Range("J89").Value = 31546.166
$ws.Range("K89").Value = 23553.3325
$ws.Range("L89").Value = 157730.83
$ws.Range("M89").Value = -17937.3325
$ws.Range("N89").Value = -168962.83

$ws.Range("H100").Value = 2475
$ws.Range("I100").Value = 2507.1428
$ws.Range("J100").Value = 2362.5
$ws.Range("K100").Value = 2507.1428
$ws.Range("L100").Value = 2362.5
$ws.Range("M100").Value = -1966.1428
$ws.Range("N100").Value = -3444.5

$ws.Range("H133").Value = 224499
$ws.Range("J133").Value = 224499
$ws.Range("L133").Value = 224499
$ws.Range("N133").Value = -234619

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 39.22222
$ws.Range("I5").Value = 38.25
$ws.Range("K5").Value = 38.25
$ws.Range("M5").Value = 73.75

$ws.Range("H61").Value = 2780.65
$ws.Range("I61").Value = 2277.9412
$ws.Range("K61").Value = 2277.9412
$ws.Range("M61").Value = -2065.9412

$ws.Range("H74").Value = 3197.9
$ws.Range("I74").Value = 3197.9
$ws.Range("K74").Value = 3197.9
$ws.Range("M74").Value = -2323.9

$ws.Range("H77").Value = 3197.9
$ws.Range("I77").Value = 3197.9
$ws.Range("K77").Value = 15989.5
$ws.Range("M77").Value = -11621.5

$ws.Range("H94").Value = 86331
$ws.Range("J94").Value = 86331
$ws.Range("L94").Value = 86331
$ws.Range("N94").Value = -88133

$ws.Range("H103").Value = 60000
$ws.Range("J103").Value = 60000
$ws.Range("L103").Value = 60000
$ws.Range("N103").Value = -62344

$ws.Range("H132").Value = 3261.25
$ws.Range("I132").Value = 2298.5715
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 6895.7145
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -4365.7145
$ws.Range("N132").Value = -35060

$ws.Range("H136").Value = 2780.65
$ws.Range("I136").Value = 2277.9412
$ws.Range("K136").Value = 6833.823600000001
$ws.Range("M136").Value = -4283.823600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 39.22222
$ws.Range("I4").Value = 38.25
$ws.Range("K4").Value = 38.25
$ws.Range("M4").Value = 76.75

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

$ws.Range("H132").Value = 100000
$ws.Range("J132").Value = 100000
$ws.Range("L132").Value = 100000
$ws.Range("N132").Value = -110120

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19024.375
$ws.Range("J41").Value = 21665.834
$ws.Range("L41").Value = 21665.834
$ws.Range("N41").Value = -22521.834

$ws.Range("H60").Value = 13281.333

$ws.Range("H132").Value = 1289.8889
$ws.Range("I132").Value = 1289.8889
$ws.Range("K132").Value = 3869.6667
$ws.Range("M132").Value = -1339.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 459
$ws.Range("J7").Value = 500
$ws.Range("L7").Value = 1500
$ws.Range("N7").Value = -1724

$ws.Range("H14").Value = 77583.164
$ws.Range("I14").Value = 77583.164
$ws.Range("K14").Value = 232749.492
$ws.Range("M14").Value = -232576.492

$ws.Range("H97").Value = 1315.1538
$ws.Range("J97").Value = 1152.125
$ws.Range("L97").Value = 3456.375
$ws.Range("N97").Value = -4448.375

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 6000
$ws.Range("M99").Value = -3754

$ws.Range("H107").Value = 1344
$ws.Range("I107").Value = 799
$ws.Range("J107").Value = 1412.125
$ws.Range("K107").Value = 2397
$ws.Range("L107").Value = 4236.375
$ws.Range("M107").Value = -477
$ws.Range("N107").Value = -8076.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H107").Value = 295.25
$ws.Range("I107").Value = 290
$ws.Range("K107").Value = 290
$ws.Range("M107").Value = 1630

$ws.Range("H122").Value = 37112.895
$ws.Range("I122").Value = 39309.465
$ws.Range("K122").Value = 117928.395
$ws.Range("M122").Value = -115478.395

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 816.63635
$ws.Range("I22").Value = 816.63635
$ws.Range("K22").Value = 816.63635
$ws.Range("M22").Value = -521.63635

$ws.Range("H27").Value = 816.63635
$ws.Range("I27").Value = 816.63635
$ws.Range("K27").Value = 816.63635
$ws.Range("M27").Value = -709.63635

$ws.Range("H63").Value = 26250
$ws.Range("I63").Value = 15000
$ws.Range("K63").Value = 15000
$ws.Range("M63").Value = -14251

$ws.Range("H66").Value = 26250
$ws.Range("I66").Value = 15000
$ws.Range("K66").Value = 45000
$ws.Range("M66").Value = -41256

$ws.Range("H132").Value = 20640.215
$ws.Range("I132").Value = 22496.3
$ws.Range("K132").Value = 67488.89999999999
$ws.Range("M132").Value = -64958.89999999999
